$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.210.77'
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").Value = '1.644.16'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''216.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = '''19.91'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").Value = '1.874.05'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").Value = '1.632.61'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("E15").Value = '  -2.77%  '
$ws.Range("D16").Value = '0.0₃0765'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '''63.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").Value = '26.220.21'
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("E20").Value = '  -0.98%  '
$ws.Range("D21").Value = '''194.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("D22").Value = '''10.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("D23").Value = '''6.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '''142.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").Value = '''6.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").Value = '''0.0503'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.89%  '
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("E34").Value = '  +1.79%  '
$ws.Range("D35").Value = '''2.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.43%  '
$ws.Range("D36").Value = '''0.911'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").Value = '1.136.23'
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("D38").Value = '''0.553'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("E39").Value = '  -1.68%  '
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '''100.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.49%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").Value = '1.782.52'
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₆0109'
$ws.Range("E46").Value = '  -2.45%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''56.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.48'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.30%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.0518'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.08%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.418'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''7.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.66%  '
